$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the old "_GoBack" bookmark (it currently sits inside the
#    "Update mesh collider (... drag model invece vuole l'oggetto)" bullet).
#    It will be re-created later at its new location.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Bold the "WF_Material" run.
# ---------------------------------------------------------------------------
$rBold = $d.Content
$null = $rBold.Find.Execute("WF_Material", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rBold.Bold = 1

# ---------------------------------------------------------------------------
# 3. Extend the "Aggiungere il materiale WF_Material..." paragraph with the
#    VertexMaterial sentence, moving the "_GoBack" bookmark to its end.
#    (Done before the paragraph-3 split below so the paragraph index for
#    paragraph 4 is still the original one.)
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs(4)
$r4 = $p4.Range
$rDot = $d.Range($r4.End - 2, $r4.End - 1)

$xmlTail = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve"> e il materiale </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:b/><w:lang w:val="it-IT"/></w:rPr><w:t>VertexMaterial</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve"> al modello (senza questa operazione i segmenti non verranno colorati cambiando il singolo triangolo).</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$rDot.InsertXML($xmlTail)

# ---------------------------------------------------------------------------
# 4. Split the "I file con le coordinate..." paragraph: simplify its mark
#    formatting, merge "cartella " + "in" into a single run, and append a
#    brand-new bullet about prefab import settings.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$r3 = $p3.Range

$xmlPara3 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:pStyle w:val="Paragrafoelenco"/>
<w:rPr>
<w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
<w:color w:val="000000"/>
<w:sz w:val="19"/>
<w:szCs w:val="19"/>
<w:lang w:val="it-IT"/>
</w:rPr>
</w:pPr>
<w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve">I file con le coordinate dei modelli, dei </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t>cage</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t>, delle coordinate baricentriche e delle annotazioni vanno inserite nella cartella in</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="it-IT"/></w:rPr><w:t>&#8220;..</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="it-IT"/></w:rPr><w:t>\</w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="it-IT"/></w:rPr><w:t>TestSelection</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="it-IT"/></w:rPr><w:t>\</w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="it-IT"/></w:rPr><w:t>Assets</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="it-IT"/></w:rPr><w:t>\</w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="it-IT"/></w:rPr><w:t>StreamingAssets</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="it-IT"/></w:rPr><w:t>&#8221;</w:t></w:r>
</w:p>
<w:p>
<w:pPr>
<w:pStyle w:val="Paragrafoelenco"/>
<w:rPr><w:b/></w:rPr>
</w:pPr>
<w:r><w:t xml:space="preserve">I prefab </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>devono</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>avere</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> come </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>opzioni</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> read/write enabled checked an</w:t></w:r>
<w:r><w:t>d optimize mesh = nothing.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$r3.InsertXML($xmlPara3)
